$d = $word.ActiveDocument

# --- Paragraph 2: "Heroku json server upload - " -----------------------
# The original text is split across several runs wrapped with w:proofErr
# spell-check markers ("Heroku" / "json"). Merge them into a single run
# and drop the now orphaned proofErr markers by round-tripping the
# paragraph's XML (Range.XML() already normalises/omits proofErr, then
# InsertXML writes that clean XML back in place of the original content).
$p2 = $d.Paragraphs(2)
$p2.Range.InsertXML($p2.Range.XML())

# --- Paragraph 4: "Blog - Objects and its internal Representation in javascript -" ---
# Merge the spell-check-split run text into a single run with the new
# wording, then clean up the leftover proofErr markers the same way.
$d.Content.Find.Execute(
    "Objects and its internal Representation in javascriptObjects and its internal Representation in javascript  - ",
    $false, $false, $false, $false, $false, $true, 1, $false,
    "Objects and its internal Representation in javascript -", 2)
$p4 = $d.Paragraphs(4)
$p4.Range.InsertXML($p4.Range.XML())

# --- Paragraph 5: "Blog - Result of running differnet values in typeof js function - " ---
$d.Content.Find.Execute(
    "Result of running differnet values in typeof js function",
    $false, $false, $false, $false, $false, $true, 1, $false,
    "Result of running differnet values in typeof js function", 2)
$p5 = $d.Paragraphs(5)
$p5.Range.InsertXML($p5.Range.XML())
